$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("atomic-properties-benchmark-6-J")
$ws.Range("C135").Value = 0.00000053524971008300805
$wb.RefreshAll()
$ws2 = $wb.Worksheets.Item("Getters")
$co = $ws2.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection().Item(4)
Write-Host "Values:" $ser.Values()
